# order_matrix_12_A_block2_VR.xlsx edit
# - participant id 12 -> 14 throughout
# - instruction video paths moved under an extra "VR" subfolder
# - rows 2..17: block_num (col I) 2 -> 4
# - video rows (3,7,11): dimension arousal -> valence; video_id / path updated
# - row 11 (col G): trailing space removed from "inverse "
# - rows 14-17 replaced with the content that used to live in rows 18-21
#   (luminance_instructions / luminance / confidence_luminance_instructions /
#   rest_suprablock), and the former rows 18-21 are removed entirely
#   (sheet shrinks from A1:L21 to A1:L17)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

function Set-NumCell($sheet, $addr, $num) {
    $sheet.Range($addr).Value = $num
}

function Clear-Cell($sheet, $addr) {
    $sheet.Range($addr).ClearContents()
}

# 1) Drop the old rows 18-21 (their content is being reused by rows 14-17 below,
#    and the sheet's used range shrinks to A1:L17).
$ws.Range("A18:L21").EntireRow.Delete()

# 2) participant id column A: 12 -> 14 for every data row (2..17)
for ($r = 2; $r -le 17; $r++) {
    Set-TextCell $ws "A$r" "14"
}

# 3) block_num column I: 2 -> 4 for rows 2..13 (rows 14-17 get no block_num, see below)
for ($r = 2; $r -le 13; $r++) {
    Set-NumCell $ws "I$r" 4
}

# Row 2 - audio_instruction (block_4_text)
Set-TextCell $ws "H2" "./instructions_videos/VR/block_4_text.mp4"

# Row 3 - video (valence, video_id 11 -> 12)
Set-TextCell $ws "E3" "12"
Set-TextCell $ws "F3" "valence"
Set-TextCell $ws "H3" "../stimuli/exp_videos/VR/12.mp4"

# Row 4 - post_stimulus_self_report
Set-TextCell $ws "H4" "./instructions_videos/VR/post_stimulus_self_report.mp4"

# Row 5 - motion_sickness
Set-TextCell $ws "H5" "./instructions_videos/VR/mareo.mp4"

# Row 6 - audio_instruction (block_4_text_reminder)
Set-TextCell $ws "H6" "./instructions_videos/VR/block_4_text_reminder.mp4"

# Row 7 - video (valence, video_id 5 -> 3)
Set-TextCell $ws "E7" "3"
Set-TextCell $ws "F7" "valence"
Set-TextCell $ws "H7" "../stimuli/exp_videos/VR/3.mp4"

# Row 8 - post_stimulus_self_report
Set-TextCell $ws "H8" "./instructions_videos/VR/post_stimulus_self_report.mp4"

# Row 9 - motion_sickness
Set-TextCell $ws "H9" "./instructions_videos/VR/mareo.mp4"

# Row 10 - audio_instruction (block_4_text_reminder)
Set-TextCell $ws "H10" "./instructions_videos/VR/block_4_text_reminder.mp4"

# Row 11 - video (valence, video_id 1 -> 2, "inverse " -> "inverse")
Set-TextCell $ws "E11" "2"
Set-TextCell $ws "F11" "valence"
Set-TextCell $ws "G11" "inverse"
Set-TextCell $ws "H11" "../stimuli/exp_videos/VR/2.mp4"

# Row 12 - post_stimulus_self_report
Set-TextCell $ws "H12" "./instructions_videos/VR/post_stimulus_self_report.mp4"

# Row 13 - motion_sickness
Set-TextCell $ws "H13" "./instructions_videos/VR/mareo.mp4"

# Row 14 - now luminance_instructions (used to be the 2nd audio_instruction reminder)
Set-TextCell $ws "H14" "./instructions_videos/VR/luminance_instructions_inverse.mp4"
Clear-Cell $ws "I14"
Set-TextCell $ws "L14" "luminance_instructions"

# Row 15 - now luminance (used to be the 4th video row)
Clear-Cell $ws "E15"
Set-TextCell $ws "F15" "luminance"
Set-TextCell $ws "H15" "../stimuli/exp_videos/VR/green_intensity_video_12.mp4"
Clear-Cell $ws "I15"
Set-TextCell $ws "L15" "luminance"

# Row 16 - now confidence_luminance_instructions (used to be post_stimulus_self_report)
Set-TextCell $ws "H16" "./instructions_videos/VR/confidence_luminance_practice_instructions_text.mp4"
Clear-Cell $ws "I16"
Set-TextCell $ws "L16" "confidence_luminance_instructions"

# Row 17 - now rest_suprablock (used to be motion_sickness)
Set-TextCell $ws "H17" "./instructions_videos/VR/rest_suprablock_text.mp4"
Clear-Cell $ws "I17"
Set-TextCell $ws "L17" "rest_suprablock"
